# "Update Data Sources from LFX" style refresh:
#  1. Six tables (slides 9, 13, 14, 15, 16, 18) switch from the default
#     table style {19266D9C-3169-4CB3-B0A4-9DD695BD7118} to
#     {A1FA7028-C1A7-44A9-B84D-1D3E6534E97C}.
#  2. The deck's theme color scheme is refreshed to the "Default" palette
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) that the regenerated
#     theme part now carries.

$p = $ppt.ActivePresentation

# --- 1. Retarget the table style on every slide that has one of the
#        affected tables. We scan shapes instead of hard-coding indices
#        so the script is resilient to any shape-order differences. ---
$targetSlides = @(9, 13, 14, 15, 16, 18)
$oldStyle = "{19266D9C-3169-4CB3-B0A4-9DD695BD7118}"
$newStyle = "{A1FA7028-C1A7-44A9-B84D-1D3E6534E97C}"

foreach ($slideNum in $targetSlides) {
    $slide = $p.Slides.Item($slideNum)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyle) {
                $tbl.ApplyStyle($newStyle)
            }
        }
    }
}

# --- 2. Update the presentation's theme color scheme to the "Default"
#        palette now stored in the (regenerated) theme part. Colors are
#        expressed as COM RGB() integers: R + G*256 + B*65536. ---
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

$themeColors.Colors(1).RGB  = 0         # dk1      #000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      #FFFFFF
$themeColors.Colors(3).RGB  = 5800213   # dk2      #158158
$themeColors.Colors(4).RGB  = 15987699  # lt2      #F3F3F3
$themeColors.Colors(5).RGB  = 13077765  # accent1  #058DC7
$themeColors.Colors(6).RGB  = 3322960   # accent2  #50B432
$themeColors.Colors(7).RGB  = 1791725   # accent3  #ED561B
$themeColors.Colors(8).RGB  = 61421     # accent4  #EDEF00
$themeColors.Colors(9).RGB  = 15059748  # accent5  #24CBE5
$themeColors.Colors(10).RGB = 7529828   # accent6  #64E572
$themeColors.Colors(11).RGB = 13369378  # hlink    #2200CC
$themeColors.Colors(12).RGB = 9116245   # folHlink #551A8B
